# Rename the leading-underscore ODM attribute column headers to their
# plain-text equivalents ("_content" -> "content", "_language" -> "language")
# across all the sheets that carry those headers in row 1.

$wb = $excel.ActiveWorkbook

# Sheet name -> cell(s) whose value needs updating.
$changes = @{
    "TranslatedText" = @{ "A1" = "language"; "C1" = "content" }
    "Title"          = @{ "A1" = "content" }
    "CheckValue"     = @{ "A1" = "content" }
    "Code"           = @{ "A1" = "content" }
    "WorkflowEnd"    = @{ "B1" = "content" }
    "UserName"       = @{ "A1" = "content" }
    "Prefix"         = @{ "A1" = "content" }
    "Suffix"         = @{ "A1" = "content" }
    "FullName"       = @{ "A1" = "content" }
    "GivenName"      = @{ "A1" = "content" }
    "FamilyName"     = @{ "A1" = "content" }
    "StreetName"     = @{ "A1" = "content" }
    "HouseNumber"    = @{ "A1" = "content" }
    "City"           = @{ "A1" = "content" }
    "StateProv"      = @{ "A1" = "content" }
    "Country"        = @{ "A1" = "content" }
    "PostalCode"     = @{ "A1" = "content" }
    "OtherText"      = @{ "A1" = "content" }
    "Meaning"        = @{ "A1" = "content" }
    "LegalReason"    = @{ "A1" = "content" }
    "DateTimeStamp"  = @{ "A1" = "content" }
    "ReasonForChange"= @{ "A1" = "content" }
    "SourceID"       = @{ "A1" = "content" }
    "FlagValue"      = @{ "B1" = "content" }
    "FlagType"       = @{ "B1" = "content" }
    "Value"          = @{ "B1" = "content" }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $changes[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
